# Update keyword-expertise scores on the "2 Chu Vanallen" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2 Chu Vanallen")

# Keyword 1 (row 2): M -> L
$ws.Range("C2").Value = "L"
# Keyword 2 (row 3): L -> (blank)
$ws.Range("C3").Value = ""
# Keyword 3 (row 4): L -> (blank)
$ws.Range("C4").Value = ""
# Keyword 6 (row 7): (blank) -> L
$ws.Range("C7").Value = "L"
# Keyword 7 (row 8): (blank) -> L
$ws.Range("C8").Value = "L"
# Keyword 8 (row 9): (blank) -> L
$ws.Range("C9").Value = "L"
# Keyword 9 (row 10): L -> (blank)
$ws.Range("C10").Value = ""
# Keyword 10 (row 11): (blank) -> L
$ws.Range("C11").Value = "L"
# Keyword 11 (row 12): (blank) -> L
$ws.Range("C12").Value = "L"

# Update the stored selection state on the sheets that were active while
# editing, matching the view state captured when the file was re-saved.
$wsCoi = $wb.Worksheets.Item("Conflicts of Interest")
$wsCoi.Activate()
$wsCoi.Range("A10:K10").Select()

$ws.Activate()
$ws.Range("A10:K10").Select()
